# Updates of local parameter files
# Rewrites the experiment grid: new model list, drops the top_p column in
# favour of a "system" column, and appends extra rows for the newer models.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row -----------------------------------------------------------
$ws.Range("E1").Value = "system"

# --- wipe the old data rows (A2:E7) so stale cells don't linger ----------
$ws.Range("A2:E7").ClearContents()

# --- new data table (model_id, model, local, temperature, system) --------
$data = @(
  @(10, "gemma2",      "'True",  0.7, "All"),
  @(15, "gemma2:27b",  "'True",  0.7, "Linux"),
  @(20, "gpt-4o-mini",  "'False", 0.7, "All"),
  @(25, "gpt-4o",       "'False", 0.7, "All"),
  @(30, "llama3",      "'True",  0.7, "All"),
  @(35, "llama3:70b",  "'True",  0.7, "Linux"),
  @(40, "phi3",        "'True",  0.7, "All"),
  @(45, "phi3:medium", "'True",  0.7, "All")
)

$row = 2
foreach ($r in $data) {
  $ws.Range("A$row").Value = $r[0]
  $ws.Range("B$row").Value = $r[1]
  $ws.Range("C$row").Value = $r[2]
  $ws.Range("D$row").Value = $r[3]
  $ws.Range("E$row").Value = $r[4]
  $row++
}

# Column C holds text that looks like booleans ("True"/"False"); the leading
# apostrophe above forces text entry. Strip the resulting quote-prefix
# formatting so the cells come out as plain shared-string text cells again.
$ws.Range("C2:C9").ClearFormats()

# --- column widths (bestFit-style sizing from the authored sheet) --------
$ws.Columns.Item(1).ColumnWidth = 7.830729166666667
$ws.Columns.Item(2).ColumnWidth = 11.166666666666666
$ws.Columns.Item(3).ColumnWidth = 4.666666666666667
$ws.Columns.Item(4).ColumnWidth = 10.498697916666666
$ws.Columns.Item(5).ColumnWidth = 6.166666666666667

# --- selection --------------------------------------------------------
$ws.Range("F10").Select() | Out-Null
